$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (old rows 6 and 7, for MuSCs -> FAPs / MuSCs -> MuSCs)
$ws.Rows("6:7").Delete()

# Row 2: FAPs / Ccl21b / Ccr10 / FAPs (target changed from ECs to FAPs)
$ws.Range("D2").Value = "FAPs"
$ws.Range("I2").Value = 0.6167434528639912
$ws.Range("J2").Value = 0.6167434528639911
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.238415
$ws.Range("N2").Value = 3.715245
$ws.Range("O2").Value = 0.4359607654144799
$ws.Range("P2").Value = 0.4359607654144798
$ws.Range("Q2").Value = 0.490989854195
$ws.Range("R2").Value = 4.418908687755
$ws.Range("S2").Value = 0.2688759477749548
$ws.Range("T2").Value = 0.2688759477749547

# Row 3: FAPs / Ccl21b / Ccr10 / MuSCs (target changed from FAPs to MuSCs)
$ws.Range("D3").Value = "MuSCs"
$ws.Range("I3").Value = 0.6167434528639912
$ws.Range("J3").Value = 0.6167434528639911
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.602242
$ws.Range("N3").Value = 4.806725999999999
$ws.Range("O3").Value = 0.5640392345855201
$ws.Range("P3").Value = 0.5640392345855201
$ws.Range("Q3").Value = 0.6352350108526665
$ws.Range("R3").Value = 5.717115097673998
$ws.Range("S3").Value = 0.3478675050890364
$ws.Range("T3").Value = 0.3478675050890364

# Row 4: MuSCs / Ccl21b / Ccr10 / FAPs (source changed FAPs->MuSCs, target changed MuSCs->FAPs)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 0.246372
$ws.Range("H4").Value = 0.739116
$ws.Range("I4").Value = 0.3832565471360088
$ws.Range("J4").Value = 0.3832565471360088
$ws.Range("M4").Value = 1.238415
$ws.Range("N4").Value = 3.715245
$ws.Range("O4").Value = 0.4359607654144799
$ws.Range("P4").Value = 0.4359607654144798
$ws.Range("Q4").Value = 0.30511078038
$ws.Range("R4").Value = 2.74599702342
$ws.Range("S4").Value = 0.1670848176395251
$ws.Range("T4").Value = 0.1670848176395251

# Row 5: MuSCs / Ccl21b / Ccr10 / MuSCs (target changed from ECs to MuSCs)
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.246372
$ws.Range("H5").Value = 0.739116
$ws.Range("I5").Value = 0.3832565471360088
$ws.Range("J5").Value = 0.3832565471360088
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.602242
$ws.Range("N5").Value = 4.806725999999999
$ws.Range("O5").Value = 0.5640392345855201
$ws.Range("P5").Value = 0.5640392345855201
$ws.Range("Q5").Value = 0.394747566024
$ws.Range("R5").Value = 3.552728094216
$ws.Range("S5").Value = 0.2161717294964837
$ws.Range("T5").Value = 0.2161717294964837
